$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("D11").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'HardwareFault']"

# Row 24
$ws.Range("D24").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E24").Value = "[]"

# Row 25
$ws.Range("D25").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E25").Value = "['Normal', 'ParamViolation']"

# Row 26
$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['SoftwareFault']"

# Row 53
$ws.Range("D53").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal', 'HardwareFault']"

# Row 54
$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E54").Value = "['SoftwareFault']"

# Row 58
$ws.Range("D58").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['Normal', 'ParamViolation']"

# Row 69
$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal']"

# Row 75
$ws.Range("D75").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E75").Value = "['Normal', 'SoftwareFault']"

# Row 83
$ws.Range("D83").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E83").Value = "['Normal', 'SurroundingEnvironment']"
